$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = "“머신 러닝 Q & AI”이 출간되었습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2025/03/30/%eb%a8%b8%ec%8b%a0-%eb%9f%ac%eb%8b%9d-q-ai%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D18").Value = "하늘에 가까운 자리"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5411"
